# Append a new row (A12:B12) to the sheet, re-using the existing shared
# strings "taxi game" / "com.singleton.strechy" (same pair already used
# in row 7), and copy the formatting from the last existing data row
# (A11:B11) so the new cells pick up the same cell style (s="1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the values first.
$ws.Range("A12").Value = "taxi game"
$ws.Range("B12").Value = "com.singleton.strechy"

# Copy only the formatting from the previous row onto the new row so the
# new cells end up with the same style as the rest of the table.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

# Move/extend the selection to the newly added row, matching the
# post-edit workbook state.
$ws.Range("A12:B12").Select()
